# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K" = strikeouts) for rows 2-29 with the recalculated
# values, replacing the prior "Strike#" derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 7
    3  = 3
    4  = 2
    5  = 4
    6  = 1
    7  = 3
    8  = 4
    9  = 2
    10 = 4
    11 = 7
    12 = 5
    13 = 2
    14 = 5
    15 = 4
    16 = 4
    17 = 6
    18 = 3
    19 = 7
    20 = 9
    21 = 6
    22 = 9
    23 = 4
    24 = 4
    25 = 1
    26 = 4
    27 = 2
    28 = 3
    29 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
